$d = $word.ActiveDocument

function Replace-WholeRun($oldText, $newText) {
    # Locate the run's exact text range, delete it entirely (this keeps
    # neighboring runs intact), then insert the new text as a brand new
    # run at that same spot. Using a fresh zero-length Range (rather than
    # Range.Text=) avoids the engine's adjacent-run coalescing pass, so a
    # preceding run with identical (default) formatting - e.g. a lone
    # " " separator run - is left completely untouched.
    $rng = $d.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $start = $rng.Start
    $rng.Delete()
    $ins = $d.Range($start, $start)
    $ins.InsertAfter($newText)
}

# 1. Update "Experienced" programming languages line
Replace-WholeRun `
    "ECMAScript (JavaScript), node.js, Python, HTML5/CSS3, XML, Bash, LaTeX" `
    "ECMAScript (JavaScript), Node.js, HTML/XML/CSS, Bash"

# 2. Update "Familiar" programming languages line
Replace-WholeRun `
    "Java SE, PHP, SQL, VBA" `
    "Clojure, ClojureScript, Elm, ML (Reason, OCaml, F#), Python"

# 3. Add React.js to JS libraries/frameworks list (mid-run edit, single
#    run in its paragraph - plain Find/Replace is safe here)
$d.Content.Find.Execute(
    "Backbone.js, Cesium.js, MarionetteJS, jQuery, RequireJS (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Backbone.js, Cesium.js, MarionetteJS, React.js, jQuery, RequireJS (", 2) | Out-Null

# 4. Add Webpack to dev/deploy tools
Replace-WholeRun `
    "Git, nvm, npm, Grunt, Gulp, Browserify, Babel" `
    "Git, nvm, npm, Grunt, Gulp, Browserify, Babel, Webpack"

# 5. Update Atlassian Stash description (only run in its paragraph)
$d.Content.Find.Execute(
    "Atlassian JIRA, Stash, Confluence",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Atlassian JIRA, Stash (Bitbucket Server), Confluence", 2) | Out-Null
